# Brands.xlsx: replace the old BrandName/Select sample data (header + 3 rows)
# with a 2-item list: Boroline (A1, keeps the bold header style) and
# Dettol (A2, now also bold). B1 becomes an empty (but still bold-styled)
# cell, and rows 3-4 are cleared out while staying part of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Wipe the previous 4-row x 2-col sample table.
$ws.Range("A1:B4").ClearContents()

# New values.
$ws.Range("A1").Value = "Boroline"
$ws.Range("A2").Value = "Dettol"

# A2 picks up the bold style that used to live on the header row.
$ws.Range("A2").Font.Bold = $true

# Keep rows 3 and 4 present (now empty) instead of disappearing entirely,
# and nudge the sheet's used range back out to B4.
$ws.Rows.Item(3).OutlineLevel = 0
$ws.Rows.Item(4).OutlineLevel = 0
$ws.Range("B4").Font.Bold = $false

# Page setup info (paper size / orientation) now present on the sheet.
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# Final selection sits on A2.
$ws.Range("A2").Select() | Out-Null
